# Apply crypto price/volume updates (GitHub Actions data refresh, 2023-10-17)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.337.41'
$ws.Range('E2').Value = '  +2.21%  '
$ws.Range('D3').Value = '1.574.58'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  +1.46%  '
$ws.Range('D5').Value = '''210.94'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').Value = '''0.489'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  +1.68%  '
$ws.Range('D8').Value = '''46.02'
$ws.Range('E8').Value = '  +4.10%  '
$ws.Range('D9').Value = '''23.72'
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('D10').Value = '''0.247'
$ws.Range('E10').Value = '  -0.78%  '
$ws.Range('D11').Value = '''0.0591'
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = '''0.0879'
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('D13').Value = '1.801.69'
$ws.Range('E13').Value = '  +0.84%  '
$ws.Range('D14').Value = '1.567.67'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '''0.523'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '''3.69'
$ws.Range('E16').Value = '  -1.15%  '
$ws.Range('D17').Value = '28.294.00'
$ws.Range('E17').Value = '  +2.20%  '
$ws.Range('D18').Value = '''62.28'
$ws.Range('E18').Value = '  -1.13%  '
$ws.Range('D19').Value = '''227.07'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('D20').Value = '''7.38'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').Value = '0.0₃0692'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').Value = '''3.93'
$ws.Range('E23').Value = '  -4.03%  '
$ws.Range('D24').Value = '''9.15'
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('D25').Value = '''1.99'
$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('D27').Value = '''14.97'
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('D28').Value = '''6.46'
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('E29').Value = '  -2.07%  '
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').Value = '''0.0463'
$ws.Range('E32').Value = '  -1.46%  '
$ws.Range('D33').Value = '''3.19'
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').Value = '''3.09'
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('D35').Value = '1.390.05'
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('D36').Value = '''1.56'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('E37').Value = '  -3.72%  '
$ws.Range('E38').Value = '  +2.33%  '
$ws.Range('E39').Value = '  +6.15%  '
$ws.Range('D40').Value = '''0.0165'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('D41').Value = '''0.530'
$ws.Range('E41').Value = '  -2.02%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = '''0.792'
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('E44').Value = '  -0.77%  '
$ws.Range('D45').Value = '''1.85'
$ws.Range('E45').Value = '  +0.43%  '
$ws.Range('D46').Value = '''0.980'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('D47').Value = '''62.15'
$ws.Range('E47').Value = '  -2.52%  '
$ws.Range('D48').Value = '1.711.48'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').Value = '''85.73'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.0519'
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₇0977'
$ws.Range('E51').Value = '  -1.24%  '
